$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of row -> @{ D = newPriceText; E = newVolumeText }
# Only columns present in the hashtable are updated for that row.
# Values are written as literal text (matching the source file's
# inlineStr cell type) by temporarily forcing a "Text" number format
# before assigning, then reverting the cell style to "Normal" so no
# residual formatting difference is left behind.
$updates = @{
    2  = @{ D = "277.06"; E = "1.37%" }
    3  = @{ D = "27.14"; E = "1.30%" }
    4  = @{ D = "4.928"; E = "0.63%" }
    5  = @{ D = "0.06411"; E = "1.44%" }
    6  = @{ D = "6.932"; E = "0.34%" }
    7  = @{ E = "-6.79%" }
    8  = @{ D = "0.8818"; E = "-0.08%" }
    9  = @{ D = "0.1521"; E = "4.11%" }
    10 = @{ D = "0.05078"; E = "-0.46%" }
    11 = @{ D = "0.07530"; E = "1.82%" }
    12 = @{ D = "0.02880"; E = "-8.19%" }
    13 = @{ D = "0.09011"; E = "-0.33%" }
    14 = @{ D = "0.001570"; E = "0.80%" }
    15 = @{ D = "0.0006412"; E = "1.40%" }
    16 = @{ D = "0.005963"; E = "-1.14%" }
    17 = @{ D = "3.457"; E = "-0.23%" }
    18 = @{ D = "3.317"; E = "-1.14%" }
    19 = @{ E = "0.01%" }
    21 = @{ D = "0.1338" }
    22 = @{ D = "3.912"; E = "0.20%" }
    23 = @{ D = "0.04428"; E = "1.92%" }
    24 = @{ D = "0.001172"; E = "-0.53%" }
    25 = @{ D = "0.003868"; E = "5.90%" }
    26 = @{ D = "0.0001201"; E = "0.08%" }
    27 = @{ E = "14.11%" }
    40 = @{ D = "0.04137"; E = "2.39%" }
    41 = @{ D = "0.006781"; E = "2.47%" }
    42 = @{ E = "1.39%" }
    43 = @{ D = "0.002252"; E = "5.72%" }
    44 = @{ D = "0.01129"; E = "-10.07%" }
    45 = @{ D = "0.00005206"; E = "-2.40%" }
    46 = @{ D = "1.487"; E = "-36.89%" }
    47 = @{ D = "0.02026"; E = "-4.39%" }
}

foreach ($row in $updates.Keys) {
    $cols = $updates[$row]
    foreach ($col in $cols.Keys) {
        $cell = $ws.Range("$col$row")
        $cell.NumberFormat = "@"
        $cell.Value = $cols[$col]
        $cell.Style = "Normal"
    }
}
